# Export review comments as XLSX
# Adds a new "All Comments" sheet, simplifies the per-tab comment sheets'
# headers (Name/Dataset/.../Author/Comment/Resolved instead of the full
# Created At / Modified At / Resolved At / Resolved By / Replies set),
# clears the now-unused "Tab" labels on the Summary sheet, and updates
# the active sheet/selection bookkeeping.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Simplify headers on the per-tab comment sheets.
#    New shared header tail for all of them: Author | Comment | Resolved
# ---------------------------------------------------------------------

$ws = $wb.Worksheets.Item("Standards")
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "Comment"
$ws.Range("D1").Value = "Resolved"
$ws.Range("E1:H1").ClearContents()
$ws.Columns.Item(1).ColumnWidth = 21.5748299319728
$ws.Columns.Item(2).ColumnWidth = 31.4319727891157
$ws.Columns.Item(3).ColumnWidth = 105.539115646259
$ws.Columns.Item(4).ColumnWidth = 14.2840136054422
$ws.Range("D1").Select()

$ws = $wb.Worksheets.Item("Datasets")
$ws.Range("A1").Value = "Dataset"
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "Comment"
$ws.Range("D1").Value = "Resolved"
$ws.Range("E1:H1").ClearContents()
$ws.Columns.Item(2).ColumnWidth = 31.4319727891157
$ws.Columns.Item(3).ColumnWidth = 105.539115646259
$ws.Columns.Item(4).ColumnWidth = 14.2840136054422
$ws.Range("C18").Select()

$ws = $wb.Worksheets.Item("Variables")
$ws.Range("A1").Value = "Dataset"
$ws.Range("B1").Value = "Variable"
$ws.Range("C1").Value = "VLM"
$ws.Range("D1").Value = "Author"
$ws.Range("E1").Value = "Comment"
$ws.Range("F1").Value = "Resolved"
$ws.Range("G1:J1").ClearContents()
$ws.Columns.Item(4).ColumnWidth = 31.4319727891157
$ws.Columns.Item(5).ColumnWidth = 105.539115646259
$ws.Columns.Item(6).ColumnWidth = 14.2840136054422
$ws.Range("F1").Select()

$ws = $wb.Worksheets.Item("Codelists")
$ws.Range("A1").Value = "Codelist"
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "Comment"
$ws.Range("D1").Value = "Resolved"
$ws.Range("E1:H1").ClearContents()
$ws.Columns.Item(2).ColumnWidth = 31.4319727891157
$ws.Columns.Item(3).ColumnWidth = 105.539115646259
$ws.Columns.Item(4).ColumnWidth = 14.2840136054422
$ws.Range("D1").Select()

$ws = $wb.Worksheets.Item("Result Displays")
$ws.Range("A1").Value = "Result Display"
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "Comment"
$ws.Range("D1").Value = "Resolved"
$ws.Range("E1:H1").ClearContents()
$ws.Columns.Item(1).ColumnWidth = 16.9880952380952
$ws.Columns.Item(2).ColumnWidth = 31.4319727891157
$ws.Columns.Item(3).ColumnWidth = 105.539115646259
$ws.Columns.Item(4).ColumnWidth = 14.2840136054422
$ws.Range("D1").Select()

$ws = $wb.Worksheets.Item("Analysis Results")
$ws.Range("A1").Value = "Result Display"
$ws.Range("B1").Value = "Analysis Result"
$ws.Range("C1").Value = "Author"
$ws.Range("D1").Value = "Comment"
$ws.Range("E1").Value = "Resolved"
$ws.Range("F1:I1").ClearContents()
$ws.Columns.Item(1).ColumnWidth = 16.9880952380952
$ws.Columns.Item(3).ColumnWidth = 31.4319727891157
$ws.Columns.Item(4).ColumnWidth = 105.539115646259
$ws.Columns.Item(5).ColumnWidth = 14.2840136054422

# ---------------------------------------------------------------------
# 2. Clear the now-unused Tab-name labels on the Summary sheet (A3:A8)
#    and move the selection there.
# ---------------------------------------------------------------------

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("A3:A8").ClearContents()
$ws.Columns.Item(1).ColumnWidth = 21.3044217687075
$ws.Columns.Item(2).ColumnWidth = 25.8962585034014
$ws.Columns.Item(3).ColumnWidth = 25.4880952380953

# ---------------------------------------------------------------------
# 3. Add the new "All Comments" sheet at the end of the workbook with
#    the raw comment-export column headers.
# ---------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newWs.Name = "All Comments"
$newWs.Range("A1").Value = "id"
$newWs.Range("B1").Value = "author"
$newWs.Range("C1").Value = "text"
$newWs.Range("D1").Value = "createdAt"
$newWs.Range("E1").Value = "modifiedAt"
$newWs.Range("F1").Value = "resolvedAt"
$newWs.Range("G1").Value = "resolvedBy"
$newWs.Range("H1").Value = "reviewCommentOids"
$newWs.Range("I1").Value = "sources"
$newWs.Columns.Item(1).ColumnWidth = 10.6870748299320
$newWs.Columns.Item(2).ColumnWidth = 11.8044217687075
$newWs.Columns.Item(3).ColumnWidth = 50.7176870748300
$newWs.Columns.Item(4).ColumnWidth = 10.6870748299320
$newWs.Columns.Item(5).ColumnWidth = 18.0544217687075
$newWs.Columns.Item(6).ColumnWidth = 17.9217687074830
$newWs.Columns.Item(7).ColumnWidth = 18.0544217687075
$newWs.Columns.Item(8).ColumnWidth = 27.7840136054422
$newWs.Columns.Item(9).ColumnWidth = 22.9268707482994
$newWs.Range("I1").Select()

# ---------------------------------------------------------------------
# 4. Make the Summary sheet the active tab/selection again.
# ---------------------------------------------------------------------

$summary = $wb.Worksheets.Item("Summary")
$summary.Activate()
$summary.Range("C15").Select()
